# Applies the "Aulas adicionadas ou atualizadas" edits to the
# "Cronograma - C07 Monitoria" schedule table.
#
# Word COM has no enumerable Run collection, so precise run-splitting
# (needed to reproduce the w:proofErr markers Word's proofing pass
# leaves around flagged words) is done via Range.InsertXML: locate the
# target paragraph text with Find, then replace it with an explicit set
# of <w:r> runs (wrapped in a minimal WordProcessingML package, the same
# shape Range.WordOpenXML itself returns).

function Set-ParaRuns {
    param(
        [string]$SearchText,
        [string]$PPrXml,
        [string]$RunsXml
    )
    $d = $word.ActiveDocument
    $r = $d.Content
    $found = $r.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $SearchText"
        return
    }
    $ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document ' + $ns + '><w:body><w:p>' + $PPrXml + $RunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 1) "Introdução a Banco de Dados..." -> split off "a" with spellStart/spellEnd
$runs1 = '<w:r><w:t xml:space="preserve">Introdução </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>a</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> Banco de Dados e Instalação do MySQL Workbench</w:t></w:r>'
Set-ParaRuns "Introdução a Banco de Dados e Instalação do MySQL Workbench" "" $runs1

# 2) "DDL (DROP, ALTER) e Comando SHOW" -> "DDL (DROP, ALTER, TRUNCATE e SHOW)"
$runs2 = '<w:r><w:t xml:space="preserve">MySQL – Comandos </w:t></w:r>' + `
    '<w:r><w:t>DDL (DROP, ALTER</w:t></w:r>' + `
    '<w:r><w:t>, TRUNCATE e SHOW</w:t></w:r>' + `
    '<w:r><w:t>)</w:t></w:r>'
Set-ParaRuns "MySQL – Comandos DDL (DROP, ALTER) e Comando SHOW" "" $runs2

# 4) (done before 3, see note below) "DCL (GRANT e REVOKE)" -> "DQL (SELECT com JOIN)"
$runs4 = '<w:r><w:t xml:space="preserve">MySQL – Comandos </w:t></w:r>' + `
    '<w:r><w:t>DQL (SELECT com JOIN)</w:t></w:r>'
Set-ParaRuns "MySQL – Comandos DCL (GRANT e REVOKE)" "" $runs4

# 3) "DQL (SELECT com JOIN)" -> "DCL (GRANT e REVOKE)"
# NOTE: run after (4) above. (4)'s replacement text equals (3)'s ORIGINAL
# search text, and (3)'s paragraph precedes (4)'s paragraph in the
# document, so Find (which always matches the first occurrence from the
# top) still lands on the right paragraph for both edits only if (4)
# runs while (3) is still untouched, and (3) runs while its own original
# text is still the first match in the document.
$runs3 = '<w:r><w:t>MySQL – Comandos</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>DCL (GRANT e REVOKE)</w:t></w:r>'
Set-ParaRuns "MySQL – Comandos DQL (SELECT com JOIN)" "" $runs3

# 5) "MySQL Stored Procedures e Functions" -> split off "e" with gramStart/gramEnd
$pPr5 = '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
$runs5 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">MySQL Stored Procedures </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>e</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> F</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>unctions</w:t></w:r>'
Set-ParaRuns "MySQL Stored Procedures e Functions" $pPr5 $runs5

# 6) "MySQL Views e Triggers" -> split off "Views" with spellStart/spellEnd
$runs6 = '<w:r><w:t xml:space="preserve">MySQL </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Views</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> e Triggers</w:t></w:r>'
Set-ParaRuns "MySQL Views e Triggers" "" $runs6
